# Rewinding to the last functional status until weathermap is integrated.
# Updates the GHI_2024-02-02 workbook: new sunrise/sunset timestamps and
# refreshed clear/cloudy-sky irradiance figures coming from the (pre-weathermap)
# production pipeline.

$wb = $excel.ActiveWorkbook

$sunrise = "2024-02-02T07:42:50"
$sunset  = "2024-02-02T17:27:06"

# ---------------------------------------------------------------------------
# Daily sheet (single summary row)
# ---------------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily")

$wsDaily.Range("E2").Value = $sunrise
$wsDaily.Range("F2").Value = $sunset

$dailyValues = @{
    "B2" = 24.724419
    "G2" = 2531.69
    "H2" = 5706.12
    "I2" = 668.33
    "J2" = 743.8099999999999
    "K2" = 3.05
    "L2" = 742.86
}
foreach ($cell in $dailyValues.Keys) {
    $wsDaily.Range($cell).Value = $dailyValues[$cell]
}

# ---------------------------------------------------------------------------
# Hourly sheet (24 hourly rows, r=2..25)
# ---------------------------------------------------------------------------
$wsHourly = $wb.Worksheets.Item("Hourly")

for ($r = 2; $r -le 25; $r++) {
    $wsHourly.Range("E$r").Value = $sunrise
    $wsHourly.Range("F$r").Value = $sunset
    # "lon" column is refreshed for every hourly row regardless of whether
    # that hour has non-zero irradiance
    $wsHourly.Range("B$r").Value = 24.724419
}

# Per-row clear/cloudy-sky values that actually changed (daylight hours only)
$hourlyValues = @{
    9  = @{ H = 0.98;  I = 9.039999999999999; J = 2.52;               K = 0.25;               L = $null; M = 0.25 }
    10 = @{ H = 72.91; I = 343.83;             J = 40.19;              K = 18.23;              L = 0;     M = 18.23 }
    11 = @{ H = 207.46; I = 591.41;            J = 66.92;              K = 51.87;              L = 0;     M = 51.87 }
    12 = @{ H = 325.66; I = 702.51;            J = 81.66;              K = 82.95999999999999;  L = 0;     M = 82.95999999999999 }
    13 = @{ H = 405.29; I = 756.2;             J = 89.64;              K = 106.83;             L = 0;     M = 106.83 }
    14 = @{ H = 435.9;  I = 773.99;            J = 92.43000000000001;  K = 115.54;             L = 0;     M = 115.54 }
    15 = @{ H = 413.9;  I = 761.39;            J = 90.43000000000001;  K = 114.53;             L = 0;     M = 114.53 }
    16 = @{ H = 341.86; I = 714.51;            J = 83.37;              K = 119.05;             L = 0.63;  M = 118.84 }
    17 = @{ H = 228.97; I = 615.72;            J = 69.95;              K = 95.68000000000001;  L = 2.41;  M = 94.94 }
    18 = @{ H = 94.48999999999999; I = 402.78; J = 46.02;              K = 37.15;              L = 0;     M = 37.15 }
    19 = @{ H = 4.27;   I = 34.74;             J = 5.2;                K = 1.72;               L = $null; M = 1.72 }
}

foreach ($r in $hourlyValues.Keys) {
    $row = $hourlyValues[$r]
    foreach ($col in $row.Keys) {
        $val = $row[$col]
        if ($null -ne $val) {
            $wsHourly.Range("$col$r").Value = $val
        }
    }
}
